$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1636
$ws.Range("F3").Value = 860
$ws.Range("F4").Value = 264
$ws.Range("F5").Value = 81
$ws.Range("F6").Value = 1174
$ws.Range("F7").Value = 791
$ws.Range("F8").Value = 823
$ws.Range("F9").Value = 1512
$ws.Range("F10").Value = 301
$ws.Range("F11").Value = 1052
$ws.Range("F13").Value = 73
$ws.Range("F16").Value = 504
$ws.Range("F17").Value = 56
$ws.Range("G17").Value = 60
$ws.Range("F18").Value = 39
$ws.Range("F22").Value = 574
$ws.Range("F23").Value = 581
$ws.Range("F24").Value = 46
$ws.Range("F27").Value = 258
$ws.Range("F28").Value = 190

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1025
$ws.Range("F5").Value = 278
$ws.Range("F9").Value = 595

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 265

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 265
$ws.Range("F3").Value = 1636
$ws.Range("F5").Value = 860
$ws.Range("F6").Value = 264
$ws.Range("F7").Value = 1025
$ws.Range("F8").Value = 81
$ws.Range("F9").Value = 1174
$ws.Range("F10").Value = 791
$ws.Range("F11").Value = 823
$ws.Range("F12").Value = 1512
$ws.Range("F13").Value = 301
$ws.Range("F14").Value = 1052
$ws.Range("F16").Value = 73
$ws.Range("F19").Value = 504
$ws.Range("F20").Value = 56
$ws.Range("G20").Value = 60
$ws.Range("F21").Value = 39
$ws.Range("F24").Value = 278
$ws.Range("F30").Value = 574
$ws.Range("F31").Value = 581
$ws.Range("F32").Value = 46
$ws.Range("F35").Value = 258
$ws.Range("F37").Value = 190
$ws.Range("F38").Value = 595
